$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

# Row 2 (Qminus1)
$ws.Range("B2").Value = 0.01545855848938353
$ws.Range("C2").Value = 0.7701062903297097
$ws.Range("D2").Value = 1.235456751497203
$ws.Range("E2").Value = 1.111511021761459
$ws.Range("F2").Value = 1.124252542560714
$ws.Range("G2").Value = 44

# Row 3 (Q0)
$ws.Range("B3").Value = 0.105076755752719
$ws.Range("C3").Value = 1.1719144532541
$ws.Range("D3").Value = 3.506314857599392
$ws.Range("E3").Value = 1.872515649493855
$ws.Range("F3").Value = 1.876375943129755
$ws.Range("G3").Value = 138

# Row 4 (Q1)
$ws.Range("B4").Value = 0.1788847452694332
$ws.Range("C4").Value = 1.315025965430722
$ws.Range("D4").Value = 7.760782161348872
$ws.Range("E4").Value = 2.785818041679835
$ws.Range("F4").Value = 2.801050726608752
$ws.Range("G4").Value = 67
